# Add a new "rating" row between the first record's "date" row and the
# second record's "SEPARATOR" row, and append a second "rating" row for
# the second record at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 (pushes the existing rows 10-18 down to 11-19)
# and fill in the new "rating" entry for the first record.
$ws.Rows("10:10").Insert()
$ws.Range("A10").Value = "rating"
$ws.Range("B10").Value = 6.66

# Append the "rating" entry for the second record as the new last row (20).
$ws.Range("A20").Value = "rating"
$ws.Range("B20").Value = 7.3

# Match the workbook's recorded selection after the edits.
$ws.Range("D19").Select()
